$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell X1: new date "09-10-2020" ---
# Force Text format first so the date-like string is stored as a literal
# string (matching the existing header cells D1:W1), not auto-converted
# to a date serial number.
$ws.Range("X1").NumberFormat = "@"
$ws.Range("X1").Value = "09-10-2020"

# Match the bold / thin-border / center-top-aligned look of the other
# header cells in row 1 (J1:W1 use this same combination).
$ws.Range("X1").Font.Bold = $true
$ws.Range("X1").Borders.Item(7).LineStyle = 1
$ws.Range("X1").Borders.Item(8).LineStyle = 1
$ws.Range("X1").Borders.Item(9).LineStyle = 1
$ws.Range("X1").Borders.Item(10).LineStyle = 1
$ws.Range("X1").HorizontalAlignment = -4108
$ws.Range("X1").VerticalAlignment = -4160

# --- Column X data values (new "09-10-2020" cumulative deceased counts) ---
$ws.Range("X2").Value = 55
$ws.Range("X3").Value = 6128
$ws.Range("X4").Value = 22
$ws.Range("X5").Value = 794
$ws.Range("X6").Value = 929
$ws.Range("X7").Value = 186
$ws.Range("X8").Value = 1158
$ws.Range("X9").Value = 2
$ws.Range("X10").Value = 5653
$ws.Range("X11").Value = 484
$ws.Range("X12").Value = 3538
$ws.Range("X13").Value = 1548
$ws.Range("X14").Value = 238
$ws.Range("X15").Value = 1291
$ws.Range("X16").Value = 775
$ws.Range("X17").Value = 9675
$ws.Range("X18").Value = 930
$ws.Range("X19").Value = 63
$ws.Range("X20").Value = 2547
$ws.Range("X21").Value = 39430
$ws.Range("X22").Value = 83
$ws.Range("X23").Value = 60
$ws.Range("X24").Value = 0
$ws.Range("X25").Value = 17
$ws.Range("X26").Value = 974
$ws.Range("X27").Value = 556
$ws.Range("X28").Value = 3741
$ws.Range("X29").Value = 1605
$ws.Range("X30").Value = 51
$ws.Range("X31").Value = 10052
$ws.Range("X32").Value = 1208
$ws.Range("X33").Value = 311
$ws.Range("X34").Value = 702
$ws.Range("X35").Value = 6245
$ws.Range("X36").Value = 5439
